$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 99
$ws1.Range("F5").Value = 2564
$ws1.Range("F6").Value = 238
$ws1.Range("F7").Value = 381

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 99
$ws4.Range("F5").Value = 2564
$ws4.Range("F6").Value = 238
$ws4.Range("F8").Value = 1
$ws4.Range("F9").Value = 381
